$d = $word.ActiveDocument

# Locate the paragraph that still has the old, single-run sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*as counsel may be heard.*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $prefix = "on {{ new_hearing_date }} at {{ new_hearing_time }} or as soon thereafter as "
    $openQuote = [char]0x201C
    $closeQuote = [char]0x201D
    $run2 = "{{ " + $openQuote
    $run3 = "counsel may be heard"
    $run4 = $closeQuote + " if users[0].attorney.there_are_any else " + $openQuote + "possible" + $closeQuote + " }}."

    $newText = $prefix + $run2 + $run3 + $run4

    $start = $target.Range.Start

    # Rewrite the paragraph's text in one shot (keeps it inside the existing run).
    $full = $d.Range($start, $target.Range.End)
    $full.Text = $newText

    # Now split that single run into four runs that line up with the template
    # pieces above (prefix / opening quote+brace / plain text / rest), by
    # nudging formatting on-and-off across each boundary so adjacent runs no
    # longer collapse back together.
    $off1 = $start + $prefix.Length
    $off2 = $off1 + $run2.Length
    $off3 = $off2 + $run3.Length
    $off4 = $start + $newText.Length

    $m1 = $d.Range($off1, $off2)
    $m1.Bold = 1
    $m1.Bold = 0

    $m2 = $d.Range($off2, $off3)
    $m2.Bold = 1
    $m2.Bold = 0

    $m3 = $d.Range($off3, $off4)
    $m3.Bold = 1
    $m3.Bold = 0

    Write-Host $d.Paragraphs(13).Range.Text
} else {
    Write-Host "Target paragraph not found"
}
